$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell {
    param($cellRef, $text)
    $rng = $ws.Range($cellRef)
    $origStyle = $rng.Style
    $rng.NumberFormat = "@"
    $rng.Value = $text
    $rng.Style = $origStyle
}

$ws.Range('D2').Value = '20.320.51'
$ws.Range('E2').Value = '  +1.99%  '

$ws.Range('D3').Value = '1.444.85'
$ws.Range('E3').Value = '  +2.70%  '

Set-TextCell 'D4' '1.004'
$ws.Range('E4').Value = '  +0.13%  '

Set-TextCell 'D5' '0.9247'
$ws.Range('E5').Value = '  -7.69%  '

Set-TextCell 'D6' '277.97'
$ws.Range('E6').Value = '  +1.90%  '

Set-TextCell 'D7' '0.3671'
$ws.Range('E7').Value = '  -0.58%  '

Set-TextCell 'D8' '0.3137'
$ws.Range('E8').Value = '  +2.45%  '

Set-TextCell 'D9' '39.27'
$ws.Range('E9').Value = '  +0.20%  '

Set-TextCell 'D10' '1.025'
$ws.Range('E10').Value = '  +3.98%  '

Set-TextCell 'D11' '0.06540'
$ws.Range('E11').Value = '  +0.28%  '

Set-TextCell 'D12' '0.9963'
$ws.Range('E12').Value = '  -0.68%  '

Set-TextCell 'D13' '5.415'
$ws.Range('E13').Value = '  +2.39%  '

Set-TextCell 'D14' '17.65'
$ws.Range('E14').Value = '  +5.08%  '

Set-TextCell 'D15' '6.101'
$ws.Range('E15').Value = '  -0.23%  '

$ws.Range('D16').Value = '1.447.51'
$ws.Range('E16').Value = '  +2.56%  '

Set-TextCell 'D17' '0.00001020'
$ws.Range('E17').Value = '  +1.84%  '

Set-TextCell 'D18' '0.9351'
$ws.Range('E18').Value = '  -6.64%  '

Set-TextCell 'D19' '0.05618'
$ws.Range('E19').Value = '  -2.24%  '

Set-TextCell 'D20' '67.49'
$ws.Range('E20').Value = '  -7.40%  '

Set-TextCell 'D21' '5.434'
$ws.Range('E21').Value = '  -2.22%  '

Set-TextCell 'D22' '14.49'
$ws.Range('E22').Value = '  +1.54%  '

Set-TextCell 'D23' '10.91'
$ws.Range('E23').Value = '  +1.53%  '

Set-TextCell 'D24' '2.270'
$ws.Range('E24').Value = '  +0.47%  '

$ws.Range('D25').Value = '20.310.80'
$ws.Range('E25').Value = '  +1.89%  '

Set-TextCell 'D26' '2.200'
$ws.Range('E26').Value = '  -0.92%  '

Set-TextCell 'D27' '135.70'
$ws.Range('E27').Value = '  -1.33%  '

Set-TextCell 'D28' '17.03'
$ws.Range('E28').Value = '  +2.15%  '

$ws.Range('D29').Value = '1.595.32'
$ws.Range('E29').Value = '  +1.75%  '

Set-TextCell 'D30' '110.56'
$ws.Range('E30').Value = '  +2.00%  '

Set-TextCell 'D31' '3.683'
$ws.Range('E31').Value = '  -3.98%  '

Set-TextCell 'D32' '0.8169'
$ws.Range('E32').Value = '  +1.56%  '

Set-TextCell 'D33' '4.886'
$ws.Range('E33').Value = '  -6.59%  '

Set-TextCell 'D34' '0.07644'
$ws.Range('E34').Value = '  -0.19%  '

Set-TextCell 'D35' '1.502'
$ws.Range('E35').Value = '  +16.93%  '

Set-TextCell 'D36' '0.05990'
$ws.Range('E36').Value = '  +3.71%  '

Set-TextCell 'D37' '4.717'
$ws.Range('E37').Value = '  -0.49%  '

Set-TextCell 'D38' '1.141'
$ws.Range('E38').Value = '  +8.35%  '

Set-TextCell 'D39' '10.31'
$ws.Range('E39').Value = '  +1.31%  '

Set-TextCell 'D40' '0.02002'
$ws.Range('E40').Value = '  -0.95%  '

Set-TextCell 'D41' '0.9354'
$ws.Range('E41').Value = '  -6.57%  '

Set-TextCell 'D42' '0.1832'
$ws.Range('E42').Value = '  -4.95%  '

Set-TextCell 'D43' '7.014'
$ws.Range('E43').Value = '  -16.26%  '

$ws.Range('B44').Value = 'TheSandbox'
$ws.Range('C44').Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
Set-TextCell 'D44' '0.5261'
$ws.Range('E44').Value = '  +0.24%  '

$ws.Range('B45').Value = 'PancakeSwap'
$ws.Range('C45').Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
Set-TextCell 'D45' '3.532'
$ws.Range('E45').Value = '  +0.53%  '

Set-TextCell 'D46' '12.02'
$ws.Range('E46').Value = '  -0.16%  '

Set-TextCell 'D47' '120.59'
$ws.Range('E47').Value = '  +9.59%  '

Set-TextCell 'D48' '0.5172'
$ws.Range('E48').Value = '  +1.95%  '

Set-TextCell 'D49' '1.774'
$ws.Range('E49').Value = '  -0.92%  '

Set-TextCell 'D50' '0.06344'
$ws.Range('E50').Value = '  +3.06%  '

Set-TextCell 'D51' '0.9895'
$ws.Range('E51').Value = '  -1.29%  '
